$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1380
$ws.Range("I43").Value = 1301
$ws.Range("J43").Value = 1399.75
$ws.Range("K43").Value = 1301
$ws.Range("L43").Value = 1399.75
$ws.Range("M43").Value = -1232
$ws.Range("N43").Value = -1537.75

$ws.Range("H62").Value = 4287524
$ws.Range("I62").Value = 6955582.5
$ws.Range("J62").Value = 18630
$ws.Range("K62").Value = 6955582.5
$ws.Range("L62").Value = 18630
$ws.Range("M62").Value = -6954958.5
$ws.Range("N62").Value = -19878

$ws.Range("H65").Value = 4287524
$ws.Range("I65").Value = 6955582.5
$ws.Range("J65").Value = 18630
$ws.Range("K65").Value = 34777912.5
$ws.Range("L65").Value = 93150
$ws.Range("M65").Value = -34774792.5
$ws.Range("N65").Value = -99390

$ws.Range("H86").Value = 1365.2858
$ws.Range("I86").Value = 715.75
$ws.Range("J86").Value = 2231.3333
$ws.Range("K86").Value = 715.75
$ws.Range("L86").Value = 2231.3333
$ws.Range("M86").Value = 407.25
$ws.Range("N86").Value = -4477.3333

$ws.Range("H89").Value = 1365.2858
$ws.Range("I89").Value = 715.75
$ws.Range("J89").Value = 2231.3333
$ws.Range("K89").Value = 3578.75
$ws.Range("L89").Value = 11156.6665
$ws.Range("M89").Value = 2037.25
$ws.Range("N89").Value = -22388.6665

$ws.Range("H106").Value = 12349013
$ws.Range("I106").Value = 13892140
$ws.Range("K106").Value = 13892140
$ws.Range("M106").Value = -13891509

$ws.Range("H112").Value = 7178218
$ws.Range("I112").Value = 680
$ws.Range("J112").Value = 8265723.5
$ws.Range("K112").Value = 2040
$ws.Range("L112").Value = 24797170.5
$ws.Range("M112").Value = -932
$ws.Range("N112").Value = -24799386.5

$ws.Range("H116").Value = 10650574
$ws.Range("I116").Value = 13844900
$ws.Range("J116").Value = 2822
$ws.Range("K116").Value = 13844900
$ws.Range("L116").Value = 2822
$ws.Range("M116").Value = -13841458
$ws.Range("N116").Value = -9706

$ws.Range("H132").Value = 173353.58
$ws.Range("I132").Value = 196747.72
$ws.Range("J132").Value = 12193.889
$ws.Range("K132").Value = 590243.16
$ws.Range("L132").Value = 36581.667
$ws.Range("M132").Value = -587713.16
$ws.Range("N132").Value = -41641.667

$ws.Range("H135").Value = 1326.1077
$ws.Range("I135").Value = 1229.6111
$ws.Range("J135").Value = 1799.8182
$ws.Range("K135").Value = 11066.4999
$ws.Range("L135").Value = 16198.3638
$ws.Range("M135").Value = -8531.499900000001
$ws.Range("N135").Value = -21268.3638

$ws.Range("H138").Value = 3867305
$ws.Range("I138").Value = 1059027.5
$ws.Range("J138").Value = 6062867.5
$ws.Range("K138").Value = 3177082.5
$ws.Range("L138").Value = 18188602.5
$ws.Range("M138").Value = -3171942.5
$ws.Range("N138").Value = -18198882.5

$ws.Range("H141").Value = 2375.8057
$ws.Range("I141").Value = 1429.4407
$ws.Range("J141").Value = 6670.846
$ws.Range("K141").Value = 4288.3221
$ws.Range("L141").Value = 20012.538
$ws.Range("M141").Value = 891.6778999999997
$ws.Range("N141").Value = -30372.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2363.8293
$ws.Range("I61").Value = 1808.2646
$ws.Range("J61").Value = 5062.2856
$ws.Range("K61").Value = 1808.2646
$ws.Range("L61").Value = 5062.2856
$ws.Range("M61").Value = -1596.2646
$ws.Range("N61").Value = -5486.2856

$ws.Range("H102").Value = 2913.5881
$ws.Range("I102").Value = 3035.4
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 3035.4
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1413.4
$ws.Range("N102").Value = -5244

$ws.Range("H132").Value = 3034.2
$ws.Range("I132").Value = 2412.182
$ws.Range("J132").Value = 4744.75
$ws.Range("K132").Value = 7236.545999999999
$ws.Range("L132").Value = 14234.25
$ws.Range("M132").Value = -4706.545999999999
$ws.Range("N132").Value = -19294.25

$ws.Range("H136").Value = 2363.8293
$ws.Range("I136").Value = 1808.2646
$ws.Range("J136").Value = 5062.2856
$ws.Range("K136").Value = 5424.793799999999
$ws.Range("L136").Value = 15186.8568
$ws.Range("M136").Value = -2874.793799999999
$ws.Range("N136").Value = -20286.8568

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1268.6666
$ws.Range("I94").Value = 1152.0834
$ws.Range("J94").Value = 1735
$ws.Range("K94").Value = 1152.0834
$ws.Range("L94").Value = 1735
$ws.Range("M94").Value = -701.0834
$ws.Range("N94").Value = -2637

$ws.Range("H134").Value = 21741562
$ws.Range("I134").Value = 25642894
$ws.Range("K134").Value = 76928682
$ws.Range("M134").Value = -76926147

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2141.6843
$ws.Range("I31").Value = 1255.4445
$ws.Range("K31").Value = 1255.4445
$ws.Range("M31").Value = -960.4445000000001

$ws.Range("H34").Value = 2141.6843
$ws.Range("I34").Value = 1255.4445
$ws.Range("K34").Value = 1255.4445
$ws.Range("M34").Value = -1053.4445

$ws.Range("H58").Value = 1138.1321
$ws.Range("I58").Value = 760.4583
$ws.Range("J58").Value = 4763.8
$ws.Range("K58").Value = 760.4583
$ws.Range("L58").Value = 4763.8
$ws.Range("M58").Value = -557.4583
$ws.Range("N58").Value = -5169.8

$ws.Range("H132").Value = 2054.6274
$ws.Range("I132").Value = 1619.4419
$ws.Range("J132").Value = 4393.75
$ws.Range("K132").Value = 4858.3257
$ws.Range("L132").Value = 13181.25
$ws.Range("M132").Value = -2328.3257
$ws.Range("N132").Value = -18241.25

$ws.Range("H134").Value = 1776.1082
$ws.Range("I134").Value = 1113.2969
$ws.Range("J134").Value = 6018.1
$ws.Range("K134").Value = 3339.8907
$ws.Range("L134").Value = 18054.3
$ws.Range("M134").Value = -804.8906999999999
$ws.Range("N134").Value = -23124.3

$ws.Range("H136").Value = 1138.1321
$ws.Range("I136").Value = 760.4583
$ws.Range("J136").Value = 4763.8
$ws.Range("K136").Value = 2281.3749
$ws.Range("L136").Value = 14291.4
$ws.Range("M136").Value = 268.6251000000002
$ws.Range("N136").Value = -19391.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1948.3462
$ws.Range("I113").Value = 1417.8462
$ws.Range("J113").Value = 2478.8462
$ws.Range("K113").Value = 1417.8462
$ws.Range("L113").Value = 2478.8462
$ws.Range("M113").Value = 752.1538
$ws.Range("N113").Value = -6818.8462

$ws.Range("H136").Value = 15709.6
$ws.Range("J136").Value = 14081.5
$ws.Range("L136").Value = 42244.5
$ws.Range("N136").Value = -47344.5

$ws.Range("H137").Value = 54999.75
$ws.Range("J137").Value = 54999.75
$ws.Range("L137").Value = 54999.75
$ws.Range("N137").Value = -65199.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 21340
$ws.Range("I22").Value = 1350
$ws.Range("J22").Value = 34666.668
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 34666.668
$ws.Range("M22").Value = -1055
$ws.Range("N22").Value = -35256.668

$ws.Range("H27").Value = 21340
$ws.Range("I27").Value = 1350
$ws.Range("J27").Value = 34666.668
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 34666.668
$ws.Range("M27").Value = -1243
$ws.Range("N27").Value = -34880.668

$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376

$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 5789.9
$ws.Range("I132").Value = 5978.5654
$ws.Range("J132").Value = 5534.647
$ws.Range("K132").Value = 17935.6962
$ws.Range("L132").Value = 16603.941
$ws.Range("M132").Value = -15405.6962
$ws.Range("N132").Value = -21663.941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = $null
$ws.Range("N44").Value = $null

$ws.Range("H136").Value = 20754.334
$ws.Range("I136").Value = 24377.166
$ws.Range("K136").Value = 73131.49800000001
$ws.Range("M136").Value = -70581.49800000001
